$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the B1:E3 block with the values that were previously held in the
# "16"/"20" trial columns (O/R and AN/AQ), effectively dropping every other
# trial column's data for this summary block.
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 439.464458638125
$ws.Range("C2").Value = 514.43989045750004
$ws.Range("D2").Value = 437.90090132812503
$ws.Range("E2").Value = 522.57770713312505

$ws.Range("B3").Value = 432.99391268812496
$ws.Range("C3").Value = 525.7632940725
$ws.Range("D3").Value = 445.52344103999997
$ws.Range("E3").Value = 522.94383455249999

# Update the selected range to match the edited block.
$ws.Range("B1:E3").Select()
